$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7947762012481689
$ws.Range("B1").Value = 2.027393817901611
$ws.Range("C1").Value = 2.668511867523193
$ws.Range("D1").Value = 2.994131088256836
$ws.Range("E1").Value = 0.8605120778083801
